$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historical GC")

# Generation Charge values (Column D) that were missing for these rows.
# Written via a scratch helper cell + TEXT() formula so the values land as
# plain text (matching the column's existing inline-string / text convention)
# instead of being auto-coerced to numbers by a direct .Value assignment.
$values = [ordered]@{
    41 = "7.3862"
    42 = "7.3862"
    43 = "5.6784"
    44 = "5.6784"
    45 = "6.1687"
    46 = "6.1687"
    50 = "7.1997"
    51 = "7.1997"
    64 = "7.5218"
    65 = "7.5218"
    70 = "7.5979"
    71 = "7.5979"
    72 = "6.5110"
    73 = "6.5110"
    81 = "6.3620"
    82 = "6.3620"
    101 = "6.3987"
    102 = "6.3987"
    103 = "6.3784"
    104 = "6.3784"
    111 = "6.3361"
    112 = "6.3361"
    113 = "6.7152"
    114 = "6.7152"
    115 = "7.1456"
    116 = "7.1456"
    118 = "7.0155"
    119 = "7.0155"
    121 = "7.0722"
    122 = "7.0722"
    125 = "6.4921"
    126 = "6.4921"
}

$helperNum = $ws.Range("Z1")
$helperTxt = $ws.Range("Z2")

foreach ($row in $values.Keys) {
    $helperNum.Value = [double]$values[$row]
    $helperTxt.Formula = '=TEXT(Z1,"0.0000")'
    $helperTxt.Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
$ws.Range("Z1:Z2").Clear()
